$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 476; this shifts existing rows 476-574 down to 477-575.
$ws.Rows("476:476").Insert()

# Populate the newly inserted row 476 with the new record's data.
$ws.Cells.Item(476, 1).Value = 8
$ws.Cells.Item(476, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(476, 3).Value = "Coquimbo"
$ws.Cells.Item(476, 4).Value = 45244
$ws.Cells.Item(476, 5).Value = 4
$ws.Cells.Item(476, 6).Value = 100112003
$ws.Cells.Item(476, 7).Value = "Ajo"
$ws.Cells.Item(476, 8).Value = "Chino"
$ws.Cells.Item(476, 9).Value = "Primera"
$ws.Cells.Item(476, 10).Value = 400
$ws.Cells.Item(476, 11).Value = 23000
$ws.Cells.Item(476, 12).Value = 24000
$ws.Cells.Item(476, 13).Value = 23500
$ws.Cells.Item(476, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(476, 15).Value = "China"
$ws.Cells.Item(476, 16).Value = 2350
$ws.Cells.Item(476, 17).Value = 10
$ws.Cells.Item(476, 18).Value = "Hortaliza"
